# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-28 with the recalculated strikeout-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 5
    3  = 6
    4  = 4
    5  = 0
    6  = 2
    7  = 3
    8  = 4
    9  = 1
    10 = 2
    11 = 2
    12 = 6
    13 = 3
    14 = 3
    15 = 8
    16 = 4
    17 = 2
    18 = 6
    19 = 2
    20 = 8
    21 = 6
    22 = 4
    23 = 3
    24 = 6
    25 = 4
    26 = 1
    27 = 2
    28 = 4
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
